# Each practice table holds 5 data rows (1, 5, 9, 13, 17) x 5 columns.
# Update the division problems in-place, cell by cell, via the Tables OM
# so each run keeps its original font/size formatting (rPr untouched).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "62÷5="  # was 98÷5=
$t.Cell(1, 2).Range.Text = "36÷8="  # was 32÷7=
$t.Cell(1, 3).Range.Text = "64÷9="  # was 98÷3=
$t.Cell(1, 4).Range.Text = "97÷8="  # was 59÷2=
$t.Cell(1, 5).Range.Text = "66÷7="  # was 12÷6=
$t.Cell(5, 1).Range.Text = "35÷9="  # was 90÷9=
$t.Cell(5, 2).Range.Text = "38÷6="  # was 52÷3=
$t.Cell(5, 3).Range.Text = "57÷5="  # was 23÷8=
$t.Cell(5, 4).Range.Text = "42÷8="  # was 39÷6=
$t.Cell(5, 5).Range.Text = "50÷4="  # was 62÷4=
$t.Cell(9, 1).Range.Text = "83÷5="  # was 24÷3=
$t.Cell(9, 2).Range.Text = "52÷7="  # was 13÷9=
$t.Cell(9, 3).Range.Text = "26÷6="  # was 84÷5=
$t.Cell(9, 4).Range.Text = "17÷3="  # was 82÷3=
$t.Cell(9, 5).Range.Text = "77÷5="  # was 36÷7=
$t.Cell(13, 1).Range.Text = "53÷6="  # was 20÷4=
$t.Cell(13, 2).Range.Text = "33÷4="  # was 13÷4=
$t.Cell(13, 3).Range.Text = "99÷5="  # was 61÷3=
$t.Cell(13, 4).Range.Text = "62÷4="  # was 33÷3=
$t.Cell(13, 5).Range.Text = "52÷4="  # was 26÷9=
$t.Cell(17, 1).Range.Text = "84÷3="  # was 35÷6=
$t.Cell(17, 2).Range.Text = "64÷2="  # was 11÷4=
$t.Cell(17, 4).Range.Text = "64÷8="  # was 72÷4=
$t.Cell(17, 5).Range.Text = "50÷5="  # was 36÷8=
